# Edit and Delete Functionality
# - Delete the member row for "Ahmed" (original row 3, id=2): rows below shift up.
# - Update the "height" value for Adnan (row 2, column F) from 5.7 to 8.2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 3 (Ahmed / Lahore / Premium / 80 / 5.5 / 2025-03-05 / Paid).
# This shifts rows 4-6 up to become rows 3-5, matching the target dimension A1:H5.
$ws.Rows(3).Delete()

# Edit the height value in row 2 (Adnan) from 5.7 to 8.2.
$ws.Range("F2").Value = 8.199999999999999
